$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.927.38"
$ws.Range("E2").Value = '  -1.66%  '

$ws.Range("D3").Value = "'2.432.26"
$ws.Range("E3").Value = '  +7.50%  '

$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = '  -0.43%  '

$ws.Range("D5").Value = "'295.30"
$ws.Range("E5").Value = '  -1.50%  '

$ws.Range("D6").Value = "'95.64"
$ws.Range("E6").Value = '  -4.20%  '

$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = '  +1.65%  '

$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").Value = "'35.15"
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").Value = "'7.10"
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("E13").Value = '  +2.39%  '

$ws.Range("D14").Value = "'2.798.90"
$ws.Range("E14").Value = '  +7.34%  '

$ws.Range("D15").Value = "'2.414.24"
$ws.Range("E15").Value = '  +6.59%  '

$ws.Range("D16").Value = "'0.844"
$ws.Range("E16").Value = '  +6.73%  '

$ws.Range("D17").Value = "'14.20"
$ws.Range("E17").Value = '  +4.42%  '

$ws.Range("D18").Value = "'45.696.39"
$ws.Range("E18").Value = '  -2.19%  '

$ws.Range("D19").Value = "'12.46"
$ws.Range("E19").Value = '  -1.62%  '

$ws.Range("D20").Value = "'0.0₃0945"
$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").Value = "'6.20"
$ws.Range("E21").Value = '  +7.14%  '

$ws.Range("D22").Value = "'67.43"
$ws.Range("E22").Value = '  +3.08%  '

$ws.Range("D23").Value = "'243.54"
$ws.Range("E23").Value = '  -1.48%  '

$ws.Range("D24").Value = "'2.79"
$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = '  +5.21%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").Value = "'39.08"
$ws.Range("E27").Value = '  -5.79%  '

$ws.Range("E28").Value = '  -0.89%  '

$ws.Range("D29").Value = "'9.75"
$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").Value = "'3.87"
$ws.Range("E30").Value = '  +17.15%  '

$ws.Range("D31").Value = "'21.31"
$ws.Range("E31").Value = '  +5.39%  '

$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("D33").Value = "'147.96"
$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("D34").Value = "'5.52"
$ws.Range("E34").Value = '  +3.14%  '

$ws.Range("D35").Value = "'0.0770"
$ws.Range("E35").Value = '  +0.50%  '

$ws.Range("D36").Value = "'2.00"
$ws.Range("E36").Value = '  +19.05%  '

$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("D39").Value = "'14.80"
$ws.Range("E39").Value = '  -5.43%  '

$ws.Range("D40").Value = "'3.81"
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("D41").Value = "'0.0300"
$ws.Range("E41").Value = '  +1.58%  '

$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = "'3.24"
$ws.Range("E42").Value = '  +4.69%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = "'1.988.45"
$ws.Range("E43").Value = '  +11.48%  '

$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("D45").Value = "'89.83"
$ws.Range("E45").Value = '  -1.11%  '

$ws.Range("E46").Value = '  -2.92%  '

$ws.Range("D47").Value = "'16.50"
$ws.Range("E47").Value = '  +29.75%  '

$ws.Range("D48").Value = "'8.60"
$ws.Range("E48").Value = '  +10.18%  '

$ws.Range("D49").Value = "'101.37"
$ws.Range("E49").Value = '  +7.89%  '

$ws.Range("D50").Value = "'2.667.98"
$ws.Range("E50").Value = '  +7.27%  '

$ws.Range("D51").Value = "'0.185"
$ws.Range("E51").Value = '  +0.21%  '
